$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
    # Update the Google Drive id for the "documentlists.xlsx" row
    $ws.Range("B2").Value = "1GfNKZ05YpkiAeWRb2OS9o_wR3lkRoeER"

    # Replace the "space rent" entry with the new document entry
    $ws.Range("A3").Value = "การเลือกซื้อและการเปลี่ยนยางรถยนต์.pdf"
    $ws.Range("B3").Value = "1lcVcRis5-qZayIRFJsadxno6jWq9YVir"

    # Remove the remaining rows (4-13) that are no longer part of the list
    $ws.Range("A4:B13").EntireRow.Delete()
}
catch {
    Write-Host "Error while editing document list: $_"
}
